$d = $word.ActiveDocument

# Update header-row heights: w:trHeight w:val="571" (28.55pt) becomes
# w:val="637" (31.85pt), in the header row of every table that has it
# (the tables whose header includes the chi-squared column).
$oldHeightPt = 571 / 20
$newHeightPt = 637 / 20
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    $r = $t.Rows.Item(1)
    if ([math]::Abs($r.Height - $oldHeightPt) -lt 0.01) {
        $r.Height = $newHeightPt
    }
}

# Fix mojibake: "χ" (chi) was corrupted to "Ï‡" (the UTF-8 bytes of χ
# mis-decoded as Windows-1252). Replace every remaining occurrence.
$tries = 0
while ($d.Content.Find.Execute("χ", $false, $false, $false, $false, $false, `
                                $true, 1, $false, "Ï‡", 2) -and $tries -lt 10) {
    $tries = $tries + 1
}

Write-Host "Row heights updated; chi replacements performed:" $tries
